# Update cryptocurrency price/volume figures per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.961.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.71%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.651.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.17%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.14%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'309.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.53%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.06%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.3896"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.21%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3832"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.08%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'51.27"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +3.67%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'1.352"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.26%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.09%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08438"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.14%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'23.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.18%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'7.093"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.76%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'7.887"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +4.05%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.00001315"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.71%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.653.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.26%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'94.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.85%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.92%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'19.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.03%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.935"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.92%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.06%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'13.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.46%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'23.955.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.63%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.462"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.90%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.978"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +5.59%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.48%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'151.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.64%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'5.424"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +2.42%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'138.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.73%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.817"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.92%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'2.491"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.09%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.833.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.75%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +6.60%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.08071"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.25%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.02960"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.06%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'6.739"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.88%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'10.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +5.40%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.2687"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.93%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.09134"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.06%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.7559"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.71%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -0.92%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.425"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.05%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'16.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.72%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.6950"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.56%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.459"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.15%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'4.090"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.58%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.000"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.06%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +0.73%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'134.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.42%  "
$ws.Range("E50").Style = "Normal"
